$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Change 1: remove the "Meta description" paragraph that follows
# the title heading.
# ---------------------------------------------------------------
$metaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description*") {
        $metaIndex = $i
        break
    }
}
if ($metaIndex -ge 1) {
    $metaPara = $d.Paragraphs.Item($metaIndex)
    [void]$metaPara.Range.Delete()
}

# ---------------------------------------------------------------
# Change 2/3: remove the trailing "Prompt: ..." paragraph (the
# image-generation prompt) and replace it with two new paragraphs:
#   - a bold "Play Ghostbusters for Free - Fun and Varied Slot Game"
#   - an italic paragraph containing the former meta-description text
# ---------------------------------------------------------------
$d2 = $word.ActiveDocument
$promptIndex = -1
for ($i = 1; $i -le $d2.Paragraphs.Count; $i++) {
    $p = $d2.Paragraphs.Item($i)
    if ($p.Range.Text -like "Prompt:*") {
        $promptIndex = $i
        break
    }
}
if ($promptIndex -lt 1) {
    $promptIndex = $d2.Paragraphs.Count
}

$promptPara = $d2.Paragraphs.Item($promptIndex)
$prevPara = $d2.Paragraphs.Item($promptIndex - 1)
[void]$promptPara.Range.Delete()

$d3 = $word.ActiveDocument
$anchor = $d3.Paragraphs.Item($promptIndex - 1)
$insertPoint = $d3.Range($anchor.Range.End, $anchor.Range.End)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Ghostbusters for Free - Fun and Varied Slot Game</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Ghostbusters, the slot game based on the cult classic movie franchise. Play for free and enjoy the fun and dynamic gameplay experience.</w:t></w:r></w:p>'
[void]$insertPoint.InsertXML($xml)
